# BIS-1002: removed "Internal Assignment" column from export.
#
# The "Internal Assignment" column (column O) is no longer exported, so its
# header and all data-row values are cleared out (formatting/styles stay).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O4:O16").ClearContents()

# Leave the selection on the range that was just cleared.
$ws.Range("O4:O16").Select()
